$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
